$wb = $excel.ActiveWorkbook

# Swap organ.text (D2) and body_part.text (A2) values on the "sample" sheet —
# they were entered the wrong way round.
$ws = $wb.Worksheets.Item("sample")

$bodyPart = $ws.Range("A2").Value()
$organ = $ws.Range("D2").Value()

$ws.Range("A2").Value = $organ
$ws.Range("D2").Value = $bodyPart

# Make the "sample" sheet the active sheet/tab, with D3 selected.
$ws.Activate()
$ws.Range("D3").Select()
